# Ticket 41 - When performing escape replacement, replace all instances.
# When replacing expressions with values, continue replacing them only one at a time.

$wb = $excel.ActiveWorkbook

$wsQuery       = $wb.Worksheets.Item("Query")
$wsPrepared    = $wb.Worksheets.Item("Prepared")
$wsLessGreater = $wb.Worksheets.Item("LessGreater")

# --- sharedStrings.xml change -------------------------------------------------
# The SQL literal embedded in the jt:forEach "items" expression on the
# LessGreater sheet (cell A2) is reformatted onto multiple lines.
$cell = $wsLessGreater.Cells.Item(2, 1)
$newText = "<jt:forEach items=""`${jdbc.execQuery('SELECT *`nFROM employee`nWHERE first_name <> \\'Randy\\'')}"" var=""employee"" >`${employee.first_name}"
$cell.Value2 = $newText

# --- styles.xml change --------------------------------------------------------
# A new cellXfs entry gets created for A2 (count goes from 9 to 10) - same
# font/fill/border as before, plus the alignment flag.
$cell.ShrinkToFit = $false

# --- sheetViews / selections ---------------------------------------------------
# Prepared: selection becomes the header row range A1:G1 instead of A3.
$wsPrepared.Range("A1:G1").Select() | Out-Null

# LessGreater: no longer the active/selected tab.
# Query: becomes the active/selected tab (and loses its old A1:G2 selection
# marker along the way).
$wsQuery.Activate() | Out-Null

Write-Host "Edit applied."
